$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("model") between vendor (C) and serial_number (D),
# shifting the remaining columns (old D..G) to (E..H).
$ws.Columns("D").Insert()

# Match the explicit column width (15) the diff sets for the new column D.
$ws.Columns("D").ColumnWidth = 14.17

# Header row
$ws.Range("D2").Value = "model"

# Row 3 (R1)
$ws.Range("D3").Value = "7206VXR"
$ws.Range("G3").Value = "5:00:00"

# Row 4 (R3)
$ws.Range("D4").Value = "OLIVE"
$ws.Range("G4").Value = "8:28:32"

# Row 5 (SW1)
$ws.Range("D5").Value = "IOSv"
$ws.Range("G5").Value = "8:26:00"
